# Auto-generated Excel COM-interop script applying numeric data updates
# to the Excalibur_Profits workbook (leve profit tracker data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 466.33334
$ws.Range("I4").Value = 466.33334
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 466.33334
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -352.33334
$ws.Range("N4").ClearContents()

$ws.Range("H28").Value = 519.4545000000001
$ws.Range("J28").Value = 1023.75
$ws.Range("L28").Value = 1023.75
$ws.Range("N28").Value = -1993.75

$ws.Range("H40").Value = 2993.4211
$ws.Range("J40").Value = 3000
$ws.Range("L40").Value = 3000
$ws.Range("N40").Value = -3350

$ws.Range("H74").Value = 7745.077
$ws.Range("I74").Value = 5841.143
$ws.Range("K74").Value = 5841.143
$ws.Range("M74").Value = -4905.143

$ws.Range("H77").Value = 7745.077
$ws.Range("I77").Value = 5841.143
$ws.Range("K77").Value = 29205.715
$ws.Range("M77").Value = -24525.715

$ws.Range("H107").Value = 237
$ws.Range("I107").Value = 223.21053
$ws.Range("J107").Value = 499
$ws.Range("K107").Value = 223.21053
$ws.Range("L107").Value = 499
$ws.Range("M107").Value = 1696.78947
$ws.Range("N107").Value = -4339

$ws.Range("H111").Value = 35716970
$ws.Range("I111").Value = 2602.3635
$ws.Range("J111").Value = 166669630
$ws.Range("K111").Value = 7807.0905
$ws.Range("L111").Value = 500008890
$ws.Range("M111").Value = -4740.0905
$ws.Range("N111").Value = -500015024

$ws.Range("H129").Value = 2728.6428
$ws.Range("J129").Value = 4999.5
$ws.Range("L129").Value = 14998.5
$ws.Range("N129").Value = -24998.5

$ws.Range("H135").Value = 35715788
$ws.Range("I135").Value = 45455996
$ws.Range("J135").Value = 1696.3334
$ws.Range("K135").Value = 409103964
$ws.Range("L135").Value = 15267.0006
$ws.Range("M135").Value = -409101429
$ws.Range("N135").Value = -20337.0006

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1618
$ws.Range("I4").Value = 300
$ws.Range("J4").Value = 2496.6667
$ws.Range("K4").Value = 300
$ws.Range("L4").Value = 2496.6667
$ws.Range("M4").Value = -184
$ws.Range("N4").Value = -2728.6667

$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H55").Value = 49999.5
$ws.Range("J55").Value = 49999.5
$ws.Range("L55").Value = 49999.5
$ws.Range("N55").Value = -50629.5

$ws.Range("H88").Value = 4758.8
$ws.Range("J88").Value = 3529.3333
$ws.Range("L88").Value = 3529.3333
$ws.Range("N88").Value = -4341.3333

$ws.Range("H91").Value = 4758.8
$ws.Range("J91").Value = 3529.3333
$ws.Range("L91").Value = 3529.3333
$ws.Range("N91").Value = -6337.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 2416.6667
$ws.Range("I19").Value = 900
$ws.Range("K19").Value = 900
$ws.Range("M19").Value = -727

$ws.Range("H22").Value = 2285.3572
$ws.Range("I22").Value = 1874.625
$ws.Range("J22").Value = 2833
$ws.Range("K22").Value = 1874.625
$ws.Range("L22").Value = 2833
$ws.Range("M22").Value = -1701.625
$ws.Range("N22").Value = -3179

$ws.Range("H43").Value = 684000
$ws.Range("J43").Value = 684000
$ws.Range("L43").Value = 684000
$ws.Range("N43").Value = -684362

$ws.Range("H86").Value = 1257.1333
$ws.Range("I86").Value = 1210.6666
$ws.Range("J86").Value = 1326.8334
$ws.Range("K86").Value = 1210.6666
$ws.Range("L86").Value = 1326.8334
$ws.Range("M86").Value = -87.66660000000002
$ws.Range("N86").Value = -3572.8334

$ws.Range("H89").Value = 1257.1333
$ws.Range("I89").Value = 1210.6666
$ws.Range("J89").Value = 1326.8334
$ws.Range("K89").Value = 6053.333000000001
$ws.Range("L89").Value = 6634.166999999999
$ws.Range("M89").Value = -437.3330000000005
$ws.Range("N89").Value = -17866.167

$ws.Range("H94").Value = 1117.2727
$ws.Range("I94").Value = 695.86664
$ws.Range("K94").Value = 695.86664
$ws.Range("M94").Value = -244.86664

$ws.Range("H107").Value = 1709.9565
$ws.Range("I107").Value = 1776.55
$ws.Range("J107").Value = 1266
$ws.Range("K107").Value = 1776.55
$ws.Range("L107").Value = 1266
$ws.Range("M107").Value = 143.45
$ws.Range("N107").Value = -5106

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1124.5
$ws.Range("I16").Value = 1124.5
$ws.Range("K16").Value = 1124.5
$ws.Range("M16").Value = -837.5

$ws.Range("H22").Value = 1098
$ws.Range("I22").Value = 797.6667
$ws.Range("J22").Value = 1999
$ws.Range("K22").Value = 797.6667
$ws.Range("L22").Value = 1999
$ws.Range("M22").Value = -447.6667
$ws.Range("N22").Value = -2699

$ws.Range("H31").Value = 8274.465
$ws.Range("I31").Value = 3249
$ws.Range("J31").Value = 8789.897000000001
$ws.Range("K31").Value = 3249
$ws.Range("L31").Value = 8789.897000000001
$ws.Range("M31").Value = -2954
$ws.Range("N31").Value = -9379.897000000001

$ws.Range("H34").Value = 8274.465
$ws.Range("I34").Value = 3249
$ws.Range("J34").Value = 8789.897000000001
$ws.Range("K34").Value = 3249
$ws.Range("L34").Value = 8789.897000000001
$ws.Range("M34").Value = -3047
$ws.Range("N34").Value = -9193.897000000001

$ws.Range("H63").Value = 85154.71000000001
$ws.Range("I63").Value = 72499.5
$ws.Range("K63").Value = 72499.5
$ws.Range("M63").Value = -71813.5

$ws.Range("H64").Value = 50000
$ws.Range("I64").Value = 50000
$ws.Range("K64").Value = 50000
$ws.Range("M64").Value = -49752

$ws.Range("H66").Value = 85154.71000000001
$ws.Range("I66").Value = 72499.5
$ws.Range("K66").Value = 217498.5
$ws.Range("M66").Value = -214066.5

$ws.Range("H67").Value = 50000
$ws.Range("I67").Value = 50000
$ws.Range("K67").Value = 50000
$ws.Range("M67").Value = -49142

$ws.Range("H69").Value = 4022.4443
$ws.Range("I69").Value = 4022.4443
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 4022.4443
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -3273.4443
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 4022.4443
$ws.Range("I72").Value = 4022.4443
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 12067.3329
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -8323.332900000001
$ws.Range("N72").ClearContents()

$ws.Range("H113").Value = 1124.5
$ws.Range("I113").Value = 1124.5
$ws.Range("K113").Value = 1124.5
$ws.Range("M113").Value = 1045.5

$ws.Range("H141").Value = 295024.5
$ws.Range("J141").Value = 295024.5
$ws.Range("L141").Value = 295024.5
$ws.Range("N141").Value = -305384.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 7976.4
$ws.Range("I108").Value = 3274
$ws.Range("K108").Value = 9822
$ws.Range("M108").Value = -6942

$ws.Range("H109").Value = 4825.385
$ws.Range("I109").Value = 2756.75
$ws.Range("K109").Value = 8270.25
$ws.Range("M109").Value = -7230.25

$ws.Range("H129").Value = 2684.1333
$ws.Range("I129").Value = 1279.875
$ws.Range("J129").Value = 4289
$ws.Range("K129").Value = 3839.625
$ws.Range("L129").Value = 12867
$ws.Range("M129").Value = 1160.375
$ws.Range("N129").Value = -22867

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 18800
$ws.Range("J32").Value = 18800
$ws.Range("L32").Value = 18800
$ws.Range("N32").Value = -19392

$ws.Range("H121").Value = 6566.3335
$ws.Range("J121").Value = 6566.3335
$ws.Range("L121").Value = 6566.3335
$ws.Range("N121").Value = -10060.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 86643.25
$ws.Range("I22").Value = 334678.66
$ws.Range("J22").Value = 3964.7778
$ws.Range("K22").Value = 334678.66
$ws.Range("L22").Value = 3964.7778
$ws.Range("M22").Value = -334383.66
$ws.Range("N22").Value = -4554.7778

$ws.Range("H27").Value = 86643.25
$ws.Range("I27").Value = 334678.66
$ws.Range("J27").Value = 3964.7778
$ws.Range("K27").Value = 334678.66
$ws.Range("L27").Value = 3964.7778
$ws.Range("M27").Value = -334571.66
$ws.Range("N27").Value = -4178.7778

$ws.Range("H46").Value = 1271.5714
$ws.Range("J46").Value = 1150
$ws.Range("L46").Value = 1150
$ws.Range("N46").Value = -1526

$ws.Range("H55").Value = 2111.2222
$ws.Range("I55").Value = 1133.6666
$ws.Range("J55").Value = 2600
$ws.Range("K55").Value = 1133.6666
$ws.Range("L55").Value = 2600
$ws.Range("M55").Value = -960.6666
$ws.Range("N55").Value = -2946

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 5007.16
$ws.Range("I107").Value = 1918.4
$ws.Range("J107").Value = 9640.299999999999
$ws.Range("K107").Value = 5755.200000000001
$ws.Range("L107").Value = 28920.9
$ws.Range("M107").Value = -3835.200000000001
$ws.Range("N107").Value = -32760.9

$ws.Range("H126").Value = 2439.6667
$ws.Range("I126").Value = 2406.4138
$ws.Range("K126").Value = 7219.241399999999
$ws.Range("M126").Value = -4749.241399999999
